$wb = $excel.ActiveWorkbook
$ws1 = $wb.Worksheets.Item("Sheet1")

$values = @(
    "Sources",
    "IP Address",
    "202.173.211.30",
    "220.77.140.189",
    "104.219.147.51",
    "118.100.228.125",
    "47.88.198.17",
    "195.204.232.90",
    "163.172.209.4",
    "181.114.23.70",
    "188.166.83.29",
    "217.182.82.44",
    "213.186.34.99",
    "45.76.171.247",
    "174.138.60.54",
    "109.234.36.6",
    "103.212.223.245",
    "185.212.254.201",
    "187.86.67.85",
    "122.248.242.239",
    "139.59.74.233",
    "64.225.116.175",
    "91.121.134.148",
    "208.113.194.28",
    "212.227.42.43",
    "222.124.154.58",
    "78.46.80.230",
    "79.137.57.130",
    "89.40.14.207",
    "190.216.224.37",
    "92.249.122.39",
    "98.158.178.141",
    "165.227.63.87",
    "36.92.117.143",
    "185.150.92.21",
    "103.27.206.205",
    "68.183.102.32",
    "83.96.252.206",
    "190.2.210.237",
    "176.9.145.194",
    "178.128.50.114",
    "1.20.97.227",
    "159.69.120.129",
    "64.150.191.231",
    "103.195.100.47",
    "185.122.58.233",
    "202.152.38.77",
    "157.230.253.101"
)

for ($i = 0; $i -lt $values.Length; $i++) {
    $ws1.Cells.Item($i + 1, 1).Value = $values[$i]
}

$ws1.Columns.Item(1).ColumnWidth = 31

$ws1.PageSetup.PaperSize = 9
$ws1.PageSetup.Orientation = 1

$ws1.Activate()
$ws1.Range("A2").Select()

Write-Host "done"
